$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation: update dSF (F) column values
$ws.Range("F10").Value = 2
$ws.Range("F12").Value = -6
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = 0
